$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shahid+257@troontechnologies.com / 12345Qwe!@#
$ws.Range("A2").Value2 = "shahid+257@troontechnologies.com"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:shahid+257@troontechnologies.com") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

$ws.Range("B2").Value2 = "12345Qwe!@#"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:12345Qwe!@#") | Out-Null
$ws.Range("B2").Style = "Hyperlink"

# Row 3: " " / '" "'  (no hyperlink, no special style)
$ws.Range("A3").Value2 = " "
$ws.Range("B3").Value2 = """ """

# Row 4: shahid+257@troontechnologies.comm / 12345Qwe!@  (same text as row 1)
$ws.Range("A4").Value2 = "shahid+257@troontechnologies.comm"
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:shahid+257@troontechnologies.comm") | Out-Null
$ws.Range("A4").Style = "Hyperlink"

$ws.Range("B4").Value2 = "12345Qwe!@"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:12345Qwe!@") | Out-Null
$ws.Range("B4").Style = "Hyperlink"

# Row 5: shahid+257@troontechnologies.com / 12345Qwe!@
$ws.Range("A5").Value2 = "shahid+257@troontechnologies.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:shahid+257@troontechnologies.com") | Out-Null
$ws.Range("A5").Style = "Hyperlink"

$ws.Range("B5").Value2 = "12345Qwe!@"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:12345Qwe!@") | Out-Null
$ws.Range("B5").Style = "Hyperlink"

# Final selection ends up on D5 (outside data range)
$ws.Range("D5").Select() | Out-Null
